$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2300.54974976347
$ws.Range("I2").Value = 161.549749763474
$ws.Range("B3").Value = 1896.35266243444
$ws.Range("I3").Value = 224.35266243444
$ws.Range("B4").Value = 3084.53702557092
$ws.Range("I4").Value = 1073.53702557092
$ws.Range("B5").Value = 3612.4772784069
$ws.Range("I5").Value = 1382.4772784069
$ws.Range("B6").Value = 4274.39575024501
$ws.Range("I6").Value = 1806.39575024501
$ws.Range("B7").Value = 5006.37135056473
$ws.Range("I7").Value = 2218.37135056473
$ws.Range("B8").Value = 5637.95274486119
$ws.Range("I8").Value = 3076.95274486119
$ws.Range("B9").Value = 5384.09103664034
$ws.Range("I9").Value = 2823.09103664034
$ws.Range("B10").Value = 6053.45986927511
$ws.Range("I10").Value = 3504.45986927511
$ws.Range("B11").Value = 4183.34369794319
$ws.Range("I11").Value = 1885.34369794319
$ws.Range("B12").Value = 3538.09146671031
$ws.Range("I12").Value = 1000.09146671031
$ws.Range("B13").Value = 3253.98636129934
$ws.Range("I13").Value = 670.98636129934
$ws.Range("B14").Value = 2472.57961449066
$ws.Range("I14").Value = 167.579614490663
$ws.Range("B15").Value = 1939.1207056887
$ws.Range("I15").Value = 66.1207056887033
$ws.Range("B16").Value = 3143.30005070225
$ws.Range("I16").Value = 258.300050702252
$ws.Range("B17").Value = 3754.01001556605
$ws.Range("I17").Value = 946.010015566054
$ws.Range("B18").Value = 4425.41264341047
$ws.Range("I18").Value = 1877.41264341047
$ws.Range("B19").Value = 5273.06741307675
$ws.Range("I19").Value = 2806.06741307675
$ws.Range("B20").Value = 5981.96522262336
$ws.Range("I20").Value = 3415.96522262336
$ws.Range("B21").Value = 5703.26229625489
$ws.Range("I21").Value = 3687.26229625489
$ws.Range("B22").Value = 6188.49374220328
$ws.Range("I22").Value = 3867.49374220328
$ws.Range("B23").Value = 4370.48533809857
$ws.Range("I23").Value = 2376.48533809857
$ws.Range("B24").Value = 3703.76486452103
$ws.Range("I24").Value = 1485.76486452103
$ws.Range("B25").Value = 3382.18231903412
$ws.Range("I25").Value = 869.182319034122
$ws.Range("B26").Value = 2591.93930963974
$ws.Range("I26").Value = 780.939309639738
$ws.Range("B27").Value = 2039.77519044255
$ws.Range("I27").Value = 478.775190442546
$ws.Range("B28").Value = 3201.17794531041
$ws.Range("I28").Value = 642.177945310406
$ws.Range("B29").Value = 3923.53189426058
$ws.Range("I29").Value = 1468.53189426058
$ws.Range("B30").Value = 4645.08510566383
$ws.Range("I30").Value = 2136.08510566383
$ws.Range("B31").Value = 5556.16554158714
$ws.Range("I31").Value = 2651.16554158714
$ws.Range("B32").Value = 6354.781755022
$ws.Range("I32").Value = 3725.781755022
$ws.Range("B33").Value = 5982.39217070932
$ws.Range("I33").Value = 3715.39217070932
$ws.Range("B34").Value = 6352.28596084523
$ws.Range("I34").Value = 4079.28596084523
$ws.Range("B35").Value = 4501.6460933452
$ws.Range("I35").Value = 2492.6460933452
$ws.Range("B36").Value = 3814.60114129959
$ws.Range("I36").Value = 2076.60114129959
$ws.Range("B37").Value = 3515.56521881807
$ws.Range("I37").Value = 1946.56521881807
$ws.Range("B38").Value = 2697.7968216493
$ws.Range("I38").Value = 1541.7968216493
$ws.Range("B39").Value = 2155.93047557963
$ws.Range("I39").Value = 197.930475579629
$ws.Range("B40").Value = 3326.57871483008
$ws.Range("I40").Value = 1118.57871483008
$ws.Range("B41").Value = 4075.41201442959
$ws.Range("I41").Value = 1865.41201442959
$ws.Range("B42").Value = 4894.875187802
$ws.Range("I42").Value = 2583.875187802
$ws.Range("B43").Value = 5907.27184807725
$ws.Range("I43").Value = 922.271848077252
$ws.Range("B44").Value = 6765.59118980005
$ws.Range("I44").Value = -6659.40881019995
$ws.Range("B45").Value = 6453.99508037672
$ws.Range("I45").Value = -6288.00491962328
$ws.Range("B46").Value = 6659.45931634619
$ws.Range("I46").Value = -118604.540683654
$ws.Range("B47").Value = 4801.42558959246
$ws.Range("I47").Value = -18309.5744104075
$ws.Range("B48").Value = 4071.67581654762
$ws.Range("I48").Value = -868.324183452378
$ws.Range("B49").Value = 3756.49938930968
$ws.Range("I49").Value = -116.500610690318
